$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151 (shifts old rows 151:192 down to 152:193)
$ws.Range("A151").EntireRow.Insert()

# Populate the new row 151 with the new weekly record
$ws.Range("A151").Value = 2
$ws.Range("B151").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C151").Value = "Coquimbo"
$ws.Range("D151").Value = 44951
$ws.Range("E151").Value = 4
$ws.Range("F151").Value = 100112043
$ws.Range("G151").Value = "Pepino ensalada"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 700
$ws.Range("K151").Value = 9000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 9500
$ws.Range("N151").Value = "$/caja 70 unidades"
$ws.Range("O151").Value = "Provincia de Limarí"
$ws.Range("P151").Value = 136
$ws.Range("Q151").Value = 70
$ws.Range("R151").Value = "Hortaliza"
